$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("C8").Value = 45
# Row 9
$ws.Range("A9").Value = 'P. point'
$ws.Range("C9").Value = 73
$ws.Range("D9").Value = "'3"
$ws.Range("E9").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F9").Value = 472
$ws.Range("G9").Value = "'34456.00"
# Row 10
$ws.Range("C10").Value = 32
$ws.Range("D10").Value = "'4"
$ws.Range("E10").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F10").Value = 662
$ws.Range("G10").Value = "'21184.00"
# Row 11
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = "'2.0"
$ws.Range("E11").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "'0.00"
# Row 12
$ws.Range("A12").Value = 'P. point'
$ws.Range("C12").Value = 91
$ws.Range("D12").Value = "'6"
$ws.Range("E12").Value = 'On board'
$ws.Range("F12").Value = 136
$ws.Range("G12").Value = "'12376.00"
# Row 13
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = "'3.0"
$ws.Range("E13").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 23
$ws.Range("G13").Value = "'207.00"
# Row 14
$ws.Range("C14").Value = 41
$ws.Range("D14").Value = "'4.0"
$ws.Range("E14").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 50
$ws.Range("G14").Value = "'2050.00"
# Row 15
$ws.Range("C15").Value = 47
$ws.Range("D15").Value = "'5.0"
$ws.Range("E15").Value = 'Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = 33
$ws.Range("G15").Value = "'1551.00"
# Row 16
$ws.Range("C16").Value = 62
$ws.Range("D16").Value = "'9.0"
$ws.Range("E16").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 219
$ws.Range("G16").Value = "'13578.00"
# Row 17
$ws.Range("A17").Value = 'Each'
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = "'10.0"
$ws.Range("E17").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F17").Value = 303
$ws.Range("G17").Value = "'3333.00"
# Row 18
$ws.Range("A18").Value = ""
$ws.Range("C18").Value = 88
$ws.Range("D18").Value = "'11.0"
$ws.Range("E18").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = "'0.00"
# Row 19
$ws.Range("A19").Value = 'R. mtr.'
$ws.Range("C19").Value = 81
$ws.Range("D19").Value = "'16"
$ws.Range("E19").Value = '20 mm'
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = "'3240.00"
# Row 20
$ws.Range("A20").Value = 'R. mtr.'
$ws.Range("C20").Value = 70
$ws.Range("D20").Value = "'17"
$ws.Range("E20").Value = '25 mm'
$ws.Range("F20").Value = 56
$ws.Range("G20").Value = "'3920.00"
# Row 21
$ws.Range("C21").Value = 61
$ws.Range("G21").Value = "'349713.00"
# Row 22
$ws.Range("A22").Value = ""
$ws.Range("C22").Value = 61
$ws.Range("D22").Value = "'14.0"
$ws.Range("E22").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = "'0.00"
# Row 23
$ws.Range("A23").Value = 'Mtr.'
$ws.Range("C23").Value = 14
$ws.Range("D23").Value = "'23"
$ws.Range("E23").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = "'280.00"
# Row 24
$ws.Range("A24").Value = ""
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = "'15.0"
$ws.Range("E24").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = "'0.00"
# Row 25
$ws.Range("A25").Value = 'Each'
$ws.Range("C25").Value = 87
$ws.Range("D25").Value = "'25"
$ws.Range("E25").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F25").Value = 1890
$ws.Range("G25").Value = "'164430.00"
# Row 26
$ws.Range("A26").Value = ""
$ws.Range("C26").Value = 61
$ws.Range("D26").Value = "'16.0"
$ws.Range("E26").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "'0.00"
# Row 27
$ws.Range("A27").Value = 'Each'
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = "'27"
$ws.Range("E27").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F27").Value = 492
$ws.Range("G27").Value = "'4920.00"
# Row 28
$ws.Range("C28").Value = 71
$ws.Range("D28").Value = "'17.0"
$ws.Range("E28").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
# Row 29
$ws.Range("C29").Value = 63
$ws.Range("G29").Value = "'11781.00"
# Row 30
$ws.Range("C30").Value = 75
$ws.Range("D30").Value = "'31"
$ws.Range("E30").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
# Row 31
$ws.Range("A31").Value = ""
$ws.Range("C31").Value = 99
$ws.Range("D31").Value = "'18.0"
$ws.Range("E31").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = "'0.00"
# Row 32
$ws.Range("A32").Value = 'Each'
$ws.Range("C32").Value = 85
$ws.Range("D32").Value = "'35"
$ws.Range("E32").Value = '8 Way (8+2)'
$ws.Range("F32").Value = 2184
$ws.Range("G32").Value = "'185640.00"
# Row 33
$ws.Range("C33").Value = 51
# Row 34
$ws.Range("C34").Value = 41
# Row 36
$ws.Range("G36").Value = "'812659.00"
$ws.Range("H36").Value = "'812659.00"
# Row 38
$ws.Range("G38").Value = "'812659.00"
$ws.Range("H38").Value = "'812659.00"
